$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price / 1h-volume snapshot (coinranking.com feed).
# Columns: A=Rank(idx), B=Coin, C=Link, D=Price, E=Volume(1h)

$ws.Range("D2").Value = '68.710.03'
$ws.Range("E2").Value = '  -0.45%  '

$ws.Range("D3").Value = '2.427.34'
$ws.Range("E3").Value = '  -1.90%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '559.04'
$ws.Range("E5").Value = '  -0.25%  '

$ws.Range("D6").Value = '161.52'
$ws.Range("E6").Value = '  -0.49%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").Value = '0.512'
$ws.Range("E8").Value = '  +1.00%  '

$ws.Range("D10").Value = '0.163'
$ws.Range("E10").Value = '  -1.63%  '

$ws.Range("E11").Value = '  -0.88%  '

$ws.Range("D12").Value = "'4.60"
$ws.Range("E12").Value = '  -5.80%  '

# Row 13: coin identity changed (rank reordering in source feed)
$ws.Range("B13").Value = 'WrappedBTC'
$ws.Range("C13").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D13").Value = '68.618.40'
$ws.Range("E13").Value = '  -0.43%  '

# Row 14: coin identity changed (rank reordering in source feed)
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").Value = "'0.0000175"
$ws.Range("E14").Value = '  +3.69%  '

$ws.Range("D15").Value = '2.876.47'
$ws.Range("E15").Value = '  -1.01%  '

$ws.Range("D16").Value = '23.12'
$ws.Range("E16").Value = '  -2.34%  '

$ws.Range("D17").Value = '2.431.01'
$ws.Range("E17").Value = '  +0.90%  '

$ws.Range("D18").Value = '10.46'
$ws.Range("E18").Value = '  -2.55%  '

$ws.Range("D19").Value = '336.42'
$ws.Range("E19").Value = '  -0.06%  '

$ws.Range("D20").Value = '6.91'
$ws.Range("E20").Value = '  -0.68%  '

$ws.Range("E21").Value = '  +0.25%  '

$ws.Range("D22").Value = '1.92'
$ws.Range("E22").Value = '  +1.71%  '

$ws.Range("D24").Value = '66.86'
$ws.Range("E24").Value = '  +0.08%  '

$ws.Range("D25").Value = '3.68'
$ws.Range("E25").Value = '  +0.38%  '

$ws.Range("D26").Value = '2.555.62'
$ws.Range("E26").Value = '  -1.72%  '

$ws.Range("E27").Value = '  +0.51%  '

$ws.Range("D28").Value = '8.16'
$ws.Range("E28").Value = '  -0.78%  '

$ws.Range("E29").Value = '  -0.56%  '

$ws.Range("D30").Value = '7.11'
$ws.Range("E30").Value = '  -1.57%  '

$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.07%  '

$ws.Range("D32").Value = '426.45'
$ws.Range("E32").Value = '  -1.61%  '

$ws.Range("E33").Value = '  +0.45%  '

$ws.Range("D35").Value = '159.68'
$ws.Range("E35").Value = '  +0.42%  '

$ws.Range("E36").Value = '  -0.10%  '

$ws.Range("D38").Value = '17.89'
$ws.Range("E38").Value = '  +0.54%  '

$ws.Range("E39").Value = '  -3.70%  '

$ws.Range("E40").Value = '  -1.42%  '

$ws.Range("E41").Value = '  +1.59%  '

$ws.Range("D42").Value = '4.33'
$ws.Range("E42").Value = '  -2.54%  '

$ws.Range("D43").Value = '1.07'
$ws.Range("E43").Value = '  -1.02%  '

$ws.Range("D44").Value = '2.03'
$ws.Range("E44").Value = '  -2.37%  '

# Row 45: coin identity changed (rank reordering in source feed)
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").Value = '3.33'
$ws.Range("E45").Value = '  -0.74%  '

# Row 46: coin identity changed (rank reordering in source feed)
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '130.96'
$ws.Range("E46").Value = '  -0.23%  '

$ws.Range("D47").Value = '0.0714'
$ws.Range("E47").Value = '  +0.18%  '

$ws.Range("D48").Value = "'0.480"
$ws.Range("E48").Value = '  -1.22%  '

$ws.Range("D49").Value = '0.555'
$ws.Range("E49").Value = '  -1.48%  '

$ws.Range("E50").Value = '  +0.57%  '

$ws.Range("E51").Value = '  +1.38%  '
